# Update "想去人数" (F column) values across sheets "展览" and "全部类型"
# and "本地生活", as published in the latest site data refresh.

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 197
$ws1.Range("F4").Value = 375
$ws1.Range("F6").Value = 846
$ws1.Range("F7").Value = 4261
$ws1.Range("F11").Value = 6257
$ws1.Range("F12").Value = 73
$ws1.Range("F14").Value = 2399
$ws1.Range("F17").Value = 498
$ws1.Range("F18").Value = 6
$ws1.Range("F19").Value = 9391
$ws1.Range("F21").Value = 2534
$ws1.Range("F23").Value = 2342
$ws1.Range("F26").Value = 251
$ws1.Range("F27").Value = 1993
$ws1.Range("F28").Value = 41
$ws1.Range("F35").Value = 97
$ws1.Range("F39").Value = 79
$ws1.Range("F40").Value = 107
$ws1.Range("F42").Value = 1580
$ws1.Range("F43").Value = 2607
$ws1.Range("F46").Value = 1260

# Sheet: 本地生活 (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 708
$ws3.Range("F3").Value = 920

# Sheet: 全部类型 (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 197
$ws4.Range("F3").Value = 708
$ws4.Range("F4").Value = 920
$ws4.Range("F6").Value = 375
$ws4.Range("F9").Value = 4261
$ws4.Range("F10").Value = 4261
$ws4.Range("F14").Value = 6257
$ws4.Range("F15").Value = 73
$ws4.Range("F16").Value = 2399
$ws4.Range("F18").Value = 498
$ws4.Range("F19").Value = 6
$ws4.Range("F20").Value = 9391
$ws4.Range("F23").Value = 2534
$ws4.Range("F24").Value = 2342
$ws4.Range("F26").Value = 251
$ws4.Range("F27").Value = 1993
$ws4.Range("F28").Value = 41
$ws4.Range("F34").Value = 97
$ws4.Range("F38").Value = 79
$ws4.Range("F39").Value = 107
$ws4.Range("F40").Value = 1580
$ws4.Range("F41").Value = 2607
$ws4.Range("F47").Value = 1260
